$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.852.76"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "'3.107.92"
$ws.Range("E3").Value = "  +3.93%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'389.80"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "'103.74"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").Value = "'37.25"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "'3.599.38"
$ws.Range("E13").Value = "  +3.87%  "
$ws.Range("D14").Value = "'18.76"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "'7.88"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "'3.094.43"
$ws.Range("E16").Value = "  +3.48%  "
$ws.Range("D17").Value = "'0.988"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "'10.86"
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("D19").Value = "'51.916.69"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "'3.20"
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "'70.01"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "'268.37"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("D26").Value = "'8.14"
$ws.Range("E26").Value = "  +3.16%  "
$ws.Range("D27").Value = "'27.17"
$ws.Range("E27").Value = "  +3.79%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'7.22"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "'10.37"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "'35.64"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").Value = "'50.45"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("D39").Value = "'0.288"
$ws.Range("E39").Value = "  +6.98%  "
$ws.Range("D40").Value = "'1.89"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").Value = "'2.61"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "'16.93"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").Value = "'129.06"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "'3.71"
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("D46").Value = "'22.33"
$ws.Range("E46").Value = "  +4.28%  "
$ws.Range("D47").Value = "'2.51"
$ws.Range("E47").Value = "  +6.63%  "
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").Value = "'2.049.99"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D50").Value = "'3.417.67"
$ws.Range("E50").Value = "  +3.89%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "'0.0327"
$ws.Range("E51").Value = "  -1.21%  "
